# Jail credit for in jail working.
# Appends 15 new case rows (790-804) to Sheet1, mirroring the existing
# 22CRB00136 / Hemmeter entries already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('22CRB00136', 'Hemmeter', 'DOMESTIC VIOLENCE', '2919.25(A)', 'No Data', 'No Contest', 'Guilty', '$ 50', '$ 25', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'ASSAULT - M1', '2903.13(A)', 'No Data', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'DOMESTIC VIOLENCE', '2919.25(A)', 'No Data', 'No Contest', 'Guilty', '$ 50', '$ 25', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'ASSAULT - M1', '2903.13(A)', 'No Data', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'DOMESTIC VIOLENCE', '2919.25(A)', 'No Data', 'No Contest', 'Guilty', '$ 50', '$ 25', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'ASSAULT - M1', '2903.13(A)', 'No Data', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'DOMESTIC VIOLENCE', '2919.25(A)', 'No Data', 'No Contest', 'Guilty', '$ 50', '$ 25', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'ASSAULT - M1', '2903.13(A)', 'No Data', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'Possession of Marijuana Drug Paraphernalia', '2925.141(C) ', 'Minor Misdemeanor', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'Driving Under Suspension FTA, Fines or Child Support', '4510.111 ', 'Unclassified Misdemeanor', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'DOMESTIC VIOLENCE', '2919.25(A)', 'No Data', 'No Contest', 'Guilty', '$ 50', '$ 25', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'ASSAULT - M1', '2903.13(A)', 'No Data', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'Possession of Marijuana Drug Paraphernalia', '2925.141(C) ', 'Minor Misdemeanor', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'Driving Under Suspension FTA, Fines or Child Support', '4510.111 ', 'Unclassified Misdemeanor', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'Traffic Control Device', '4511.12', 'Minor Misdemeanor', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None')
)

$startRow = 790
$r = $startRow
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $cell = $ws.Cells.Item($r, $c)
        # Force text storage so values such as "$ 50", "None", "4511.12"
        # and codes with trailing spaces round-trip as literal strings
        # instead of being coerced to numbers/currency by Excel.
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
